$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.531.57'
$ws.Range("D2").Style = $ws.Range("B2").Style
$ws.Range("E2").Value = '  -2.39%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.005.58'
$ws.Range("D3").Style = $ws.Range("B3").Style
$ws.Range("E3").Value = '  -0.41%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.18'
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = '  -8.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.600'
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = '  -3.12%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.85'
$ws.Range("D8").Style = $ws.Range("B8").Style
$ws.Range("E8").Value = '  -3.51%  '
$ws.Range("E9").Value = '  -3.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.35'
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = '  +3.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0748'
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = '  -3.56%  '
$ws.Range("E12").Value = '  -3.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.298.29'
$ws.Range("D13").Style = $ws.Range("B13").Style
$ws.Range("E13").Value = '  -0.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.22'
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E14").Value = '  -0.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.26'
$ws.Range("D15").Style = $ws.Range("B15").Style
$ws.Range("E15").Value = '  -4.35%  '
$ws.Range("E16").Value = '  -5.49%  '
$ws.Range("E17").Value = '  -2.75%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.008.76'
$ws.Range("D18").Style = $ws.Range("B18").Style
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '36.502.33'
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Value = '  -2.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '67.83'
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").Value = '  -3.24%  '
$ws.Range("E21").Value = '  -4.06%  '
$ws.Range("E22").Value = '  +3.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '221.77'
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value = '  -4.54%  '
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.38'
$ws.Range("D25").Style = $ws.Range("B25").Style
$ws.Range("E25").Value = '  +1.52%  '
$ws.Range("E26").Value = '  -8.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.94'
$ws.Range("D27").Style = $ws.Range("B27").Style
$ws.Range("E27").Value = '  -1.72%  '
$ws.Range("E28").Value = '  -3.28%  '
$ws.Range("E29").Value = '  -1.68%  '
$ws.Range("E30").Value = '  +2.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '18.82'
$ws.Range("D31").Style = $ws.Range("B31").Style
$ws.Range("E31").Value = '  -4.25%  '
$ws.Range("E32").Value = '  -2.93%  '
$ws.Range("E33").Value = '  -4.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0604'
$ws.Range("D34").Style = $ws.Range("B34").Style
$ws.Range("E34").Value = '  -6.43%  '
$ws.Range("E35").Value = '  +3.22%  '
$ws.Range("E36").Value = '  -6.21%  '
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.37'
$ws.Range("D38").Style = $ws.Range("B38").Style
$ws.Range("E38").Value = '  +0.90%  '
$ws.Range("E39").Value = '  -2.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.72'
$ws.Range("D40").Style = $ws.Range("B40").Style
$ws.Range("E40").Value = '  +6.30%  '
$ws.Range("E41").Value = '  -1.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.451.27'
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E42").Value = '  +2.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0924'
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Value = '  -0.43%  '
$ws.Range("E44").Value = '  -4.70%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.31'
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("E46").Value = '  -8.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.23'
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").Value = '  -3.83%  '
$ws.Range("E48").Value = '  -2.50%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.82'
$ws.Range("D49").Style = $ws.Range("B49").Style
$ws.Range("E49").Value = '  +24.26%  '
$ws.Range("E50").Value = '  -1.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.83'
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Value = '  -3.09%  '
